$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drohne")

# --- Step 1: fill in the product names for the three new rows (35-37) ---
# (text entered first, in row order, matches the shared-string insertion order)
$ws.Range("A35").Value = "GPS 2"
$ws.Range("C35").Value = "NEO-6M"

$ws.Range("A36").Value = "Ultraschall"

$ws.Range("A37").Value = "Höhensensor"
$ws.Range("C37").Value = "BMP180"

# --- Step 2: quantities ---
$ws.Range("B35").Value = 1
$ws.Range("B36").Value = 1
$ws.Range("B37").Value = 1

# --- Step 3: prices ---
$ws.Range("D35").Value = 11
$ws.Range("D36").Value = 2.85
$ws.Range("D37").Value = 3.89

# --- Step 4: purchase dates ---
$ws.Range("H35").Value = 43081
$ws.Range("H36").Value = 43081
$ws.Range("H37").Value = 43080

# --- Step 5: weblinks (I column) - added afterwards, including retroactively on row 34 ---
$ws.Range("I34").Value = "https://www.ebay.de/itm/252715059842"
$ws.Range("I35").Value = "https://www.ebay.de/itm/252715059842"
$ws.Range("I36").Value = "https://www.ebay.de/itm/162571876188"
$ws.Range("I37").Value = "https://www.ebay.de/itm/162675931015"

# --- Step 6: formatting to match the surrounding rows ---
# Column A: "Gut" (green) highlight style used by the other sensor rows
$ws.Range("A35").Style = "Gut"
$ws.Range("A36").Style = "Gut"
$ws.Range("A37").Style = "Gut"

# Column B: centered quantity cell (copy formatting from row above)
$ws.Range("B34").Copy($ws.Range("B35"))
$ws.Range("B34").Copy($ws.Range("B36"))
$ws.Range("B34").Copy($ws.Range("B37"))
$ws.Range("B35").Value = 1
$ws.Range("B36").Value = 1
$ws.Range("B37").Value = 1

# Column D: currency number format, matching the other price cells
$ws.Range("D30").Copy($ws.Range("D35"))
$ws.Range("D30").Copy($ws.Range("D36"))
$ws.Range("D30").Copy($ws.Range("D37"))
$ws.Range("D35").Value = 11
$ws.Range("D36").Value = 2.85
$ws.Range("D37").Value = 3.89

# Column H: date number format, matching the other date cells
$ws.Range("H30").Copy($ws.Range("H35"))
$ws.Range("H30").Copy($ws.Range("H36"))
$ws.Range("H30").Copy($ws.Range("H37"))
$ws.Range("H35").Value = 43081
$ws.Range("H36").Value = 43081
$ws.Range("H37").Value = 43080

# --- Step 7: view/selection bookkeeping to match the saved workbook state ---
$ws.Range("I37").Select()

# --- Step 8: force a full recalculation so cached formula values (e.g. the
#     Preis sum in D1) reflect the newly added rows before the workbook saves ---
$excel.Calculate()
